# Daily attendance processing - 2026-01-24 03:29:41
# Normalize the "Recorded By" (column G) values: swap the order of the
# "dnasr281@gmail.com, System" entries to "System, dnasr281@gmail.com"
# so the System actor is listed first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
